$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 15395.866
$ws.Range("I15").Value = 15395.866
$ws.Range("K15").Value = 46187.598
$ws.Range("M15").Value = -46018.598
$ws.Range("H88").Value = 1641
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 1801.25
$ws.Range("K88").Value = 1000
$ws.Range("L88").Value = 1801.25
$ws.Range("M88").Value = -594
$ws.Range("N88").Value = -2613.25
$ws.Range("H91").Value = 1641
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 1801.25
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 1801.25
$ws.Range("M91").Value = 404
$ws.Range("N91").Value = -4609.25
$ws.Range("H98").Value = 30935722
$ws.Range("I98").Value = 10527493
$ws.Range("J98").Value = 127874820
$ws.Range("K98").Value = 10527493
$ws.Range("L98").Value = 127874820
$ws.Range("M98").Value = -10525995
$ws.Range("N98").Value = -127877816
$ws.Range("H112").Value = 2138
$ws.Range("I112").Value = 700
$ws.Range("J112").Value = 2497.5
$ws.Range("K112").Value = 2100
$ws.Range("L112").Value = 7492.5
$ws.Range("M112").Value = -992
$ws.Range("N112").Value = -9708.5
$ws.Range("H116").Value = 13741005
$ws.Range("I116").Value = 9526938
$ws.Range("J116").Value = 16690852
$ws.Range("K116").Value = 9526938
$ws.Range("L116").Value = 16690852
$ws.Range("M116").Value = -9523496
$ws.Range("N116").Value = -16697736
$ws.Range("H122").Value = 30935722
$ws.Range("I122").Value = 10527493
$ws.Range("J122").Value = 127874820
$ws.Range("K122").Value = 31582479
$ws.Range("L122").Value = 383624460
$ws.Range("M122").Value = -31580029
$ws.Range("N122").Value = -383629360
$ws.Range("H132").Value = 3705083
$ws.Range("I132").Value = 1313.8889
$ws.Range("J132").Value = 37039004
$ws.Range("K132").Value = 3941.6667
$ws.Range("L132").Value = 111117012
$ws.Range("M132").Value = -1411.6667
$ws.Range("N132").Value = -111122072
$ws.Range("H137").Value = 15789430
$ws.Range("I137").Value = 15625831
$ws.Range("J137").Value = 16101048
$ws.Range("K137").Value = 46877493
$ws.Range("L137").Value = 48303144
$ws.Range("M137").Value = -46874943
$ws.Range("N137").Value = -48308244

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4383.97
$ws.Range("I32").Value = 4175.736
$ws.Range("K32").Value = 4175.736
$ws.Range("M32").Value = -3888.736
$ws.Range("H45").Value = 278673.1
$ws.Range("I45").Value = 417296.5
$ws.Range("J45").Value = 1426.25
$ws.Range("K45").Value = 417296.5
$ws.Range("L45").Value = 1426.25
$ws.Range("M45").Value = -416919.5
$ws.Range("N45").Value = -2180.25
$ws.Range("H74").Value = 36961628
$ws.Range("I74").Value = 37713250
$ws.Range("J74").Value = 35558596
$ws.Range("K74").Value = 37713250
$ws.Range("L74").Value = 35558596
$ws.Range("M74").Value = -37712376
$ws.Range("N74").Value = -35560344
$ws.Range("H77").Value = 36961628
$ws.Range("I77").Value = 37713250
$ws.Range("J77").Value = 35558596
$ws.Range("K77").Value = 188566250
$ws.Range("L77").Value = 177792980
$ws.Range("M77").Value = -188561882
$ws.Range("N77").Value = -177801716
$ws.Range("H122").Value = 3066.6667
$ws.Range("I122").Value = 2350
$ws.Range("J122").Value = 3640
$ws.Range("K122").Value = 7050
$ws.Range("L122").Value = 10920
$ws.Range("M122").Value = -4600
$ws.Range("N122").Value = -15820
$ws.Range("H132").Value = 18690856
$ws.Range("I132").Value = 14449084
$ws.Range("J132").Value = 31946394
$ws.Range("K132").Value = 43347252
$ws.Range("L132").Value = 95839182
$ws.Range("M132").Value = -43344722
$ws.Range("N132").Value = -95844242

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1208.6842
$ws.Range("I99").Value = 856.7857
$ws.Range("J99").Value = 2194
$ws.Range("K99").Value = 856.7857
$ws.Range("L99").Value = 2194
$ws.Range("M99").Value = 641.2143
$ws.Range("N99").Value = -5190
$ws.Range("H134").Value = 25211286
$ws.Range("I134").Value = 33334514
$ws.Range("J134").Value = 4903213.5
$ws.Range("K134").Value = 100003542
$ws.Range("L134").Value = 14709640.5
$ws.Range("M134").Value = -100001007
$ws.Range("N134").Value = -14714710.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 13000
$ws.Range("J55").Value = 13000
$ws.Range("L55").Value = 13000
$ws.Range("N55").Value = -13630
$ws.Range("H99").Value = 11920780
$ws.Range("I99").Value = 18528190
$ws.Range("J99").Value = 27442.8
$ws.Range("K99").Value = 18528190
$ws.Range("L99").Value = 27442.8
$ws.Range("M99").Value = -18526692
$ws.Range("N99").Value = -30438.8
$ws.Range("H126").Value = 11920780
$ws.Range("I126").Value = 18528190
$ws.Range("J126").Value = 27442.8
$ws.Range("K126").Value = 55584570
$ws.Range("L126").Value = 82328.39999999999
$ws.Range("M126").Value = -55582100
$ws.Range("N126").Value = -87268.39999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 503355.75
$ws.Range("J107").Value = 1318.0555
$ws.Range("L107").Value = 3954.1665
$ws.Range("N107").Value = -7794.166499999999
$ws.Range("H117").Value = 1902.7693
$ws.Range("J117").Value = 1902.7693
$ws.Range("L117").Value = 5708.3079
$ws.Range("N117").Value = -12592.3079

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 800
$ws.Range("I41").Value = 800
$ws.Range("K41").Value = 800
$ws.Range("M41").Value = -445
$ws.Range("H102").Value = 3166.617
$ws.Range("I102").Value = 3562.5715
$ws.Range("J102").Value = 2011.75
$ws.Range("K102").Value = 3562.5715
$ws.Range("L102").Value = 2011.75
$ws.Range("M102").Value = -1940.5715
$ws.Range("N102").Value = -5255.75
$ws.Range("H122").Value = 7967487
$ws.Range("I122").Value = 39846.5
$ws.Range("J122").Value = 33335936
$ws.Range("K122").Value = 119539.5
$ws.Range("L122").Value = 100007808
$ws.Range("M122").Value = -117089.5
$ws.Range("N122").Value = -100012708

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1440.1177
$ws.Range("I7").Value = 1583.2307
$ws.Range("J7").Value = 975
$ws.Range("K7").Value = 1583.2307
$ws.Range("L7").Value = 975
$ws.Range("M7").Value = -1471.2307
$ws.Range("N7").Value = -1199
$ws.Range("H40").Value = 3475843
$ws.Range("I40").Value = 4833990
$ws.Range("J40").Value = 5022.5557
$ws.Range("K40").Value = 4833990
$ws.Range("L40").Value = 5022.5557
$ws.Range("M40").Value = -4833854
$ws.Range("N40").Value = -5294.5557
$ws.Range("H126").Value = 1440.1177
$ws.Range("I126").Value = 1583.2307
$ws.Range("J126").Value = 975
$ws.Range("K126").Value = 4749.6921
$ws.Range("L126").Value = 2925
$ws.Range("M126").Value = -2279.6921
$ws.Range("N126").Value = -7865
$ws.Range("H132").Value = 2202818.8
$ws.Range("I132").Value = 3350495.5
$ws.Range("J132").Value = 3104.5833
$ws.Range("K132").Value = 10051486.5
$ws.Range("L132").Value = 9313.749899999999
$ws.Range("M132").Value = -10048956.5
$ws.Range("N132").Value = -14373.7499
$ws.Range("H136").Value = 1471208.6
$ws.Range("I136").Value = 1783103
$ws.Range("J136").Value = 849.2857
$ws.Range("K136").Value = 5349309
$ws.Range("L136").Value = 2547.8571
$ws.Range("M136").Value = -5346759
$ws.Range("N136").Value = -7647.8571

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5557087
$ws.Range("I122").Value = 13890341
$ws.Range("J122").Value = 1584.1666
$ws.Range("K122").Value = 41671023
$ws.Range("L122").Value = 4752.4998
$ws.Range("M122").Value = -41668573
$ws.Range("N122").Value = -9652.4998
$ws.Range("H126").Value = 18708194
$ws.Range("I126").Value = 20677388
$ws.Range("J126").Value = 848.5
$ws.Range("K126").Value = 62032164
$ws.Range("L126").Value = 2545.5
$ws.Range("M126").Value = -62029694
$ws.Range("N126").Value = -7485.5
